# Actualizacion Datos Personales 4 nov
# Applies the "Nov 4" data update to both roster sheets (3APV / 5APV):
#  - A few huge phone numbers that were stored as text got re-typed and
#    Excel re-rendered them in (lower-case) scientific notation text.
#  - Tutor e-mails / names were added, corrected, or (in one case) cleared.

$wb = $excel.ActiveWorkbook

function Set-TextValue($cell, [string]$text) {
    # Force the literal text into the cell (Excel would otherwise coerce a
    # numeric-looking string like "5.22722e+16" into a real number), then
    # drop the quote-prefix styling so the cell keeps its original
    # (default/"Normal") style - only the shared-string content should change.
    $cell.Value = "'" + $text
    $cell.Style = "Normal"
}

# ---------------------------------------------------------------------
# Sheet "3APV"
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("3APV")

# Tel_Fijo (col G) values that now render as scientific notation text
Set-TextValue $ws1.Cells.Item(5, 7) "5.22722e+16"
Set-TextValue $ws1.Cells.Item(16, 7) "2.72721e+16"

# Correo_Tutor (col I) additions / correction
$ws1.Cells.Item(3, 9).Value = "Marchate1986@gmail.com"
$ws1.Cells.Item(4, 9).Value = "jdcastro@gmail.com"
$ws1.Cells.Item(20, 9).Value = "oficialsanchezmz@gmail.com"

# ---------------------------------------------------------------------
# Sheet "5APV"
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("5APV")

# Telefono_Tutor (col J) value that now renders as scientific notation text
Set-TextValue $ws2.Cells.Item(18, 10) "5.26145e+16"

# Correo (col E) correction
$ws2.Cells.Item(38, 5).Value = "hersonxolio@gmail.com"

# Tutor (col H) additions
$ws2.Cells.Item(27, 8).Value = "BARSIMEO ISMAEL PAZ CESAR"
$ws2.Cells.Item(38, 8).Value = "ALVARO ELIAS XOLIO LINARES"

# Telefono_Tutor (col J) cleared
$ws2.Cells.Item(6, 10).ClearContents()
